$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B numeric values and D2
$ws.Range("B2").Value = 27
$ws.Range("D2").Value = 100

$ws.Range("B3").Value = 25

$ws.Range("B4").Value = 21

$ws.Range("B5").Value = 27

$ws.Range("B6").Value = 50

# Rename basetech_EL_* labels to basetech_EH_* and update values
$ws.Range("A7").Value = "basetech_EH_electricity"
$ws.Range("B7").Value = 13

$ws.Range("A8").Value = "basetech_EH_oil"
$ws.Range("B8").Value = 12

$ws.Range("A9").Value = "basetech_EH_gas"
$ws.Range("B9").Value = 11

$ws.Range("A10").Value = "basetech_EH_K"
$ws.Range("B10").Value = 13

# Update the active selection to B5
$ws.Range("B5").Select()

$wb.Save()
